# Weekly data refresh: insert the newest "Albahaca" price record for
# "Feria Lagunitas de Puerto Montt" at the top of this week's data block
# (row 88), pushing all older records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 88 (existing rows 88..207 shift to 89..208).
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row with the latest observation.
$ws.Range("A88").Value = 4
$ws.Range("B88").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C88").Value = "Los Lagos"
$ws.Range("D88").Value = 45174
$ws.Range("D88").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E88").Value = 10
$ws.Range("F88").Value = 100112052
$ws.Range("G88").Value = "Albahaca"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 90
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 6000
$ws.Range("M88").Value = 6000
$ws.Range("N88").Value = "$/paquete"
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 6000
$ws.Range("Q88").Value = 1
$ws.Range("R88").Value = "Hortaliza"
